$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Is the New Visa Bitcoin Rewards Card Worth It?", "https://twocents.lifehacker.com/is-the-new-visa-bitcoin-rewards-card-worth-it-1845803159"),
    @("Bitcoin passes `$20k and reaches all-time high", "http://techcrunch.com/2020/12/16/bitcoin-passes-20k-and-reaches-all-time-high/"),
    @("Bitcoin Blows Past `$20,000 Milestone", "https://gizmodo.com/bitcoin-blows-past-20-000-milestone-1845894176"),
    @("Is rising usage driving crypto’s recent price boom?", "http://techcrunch.com/2020/12/17/is-rising-usage-driving-cryptos-recent-price-boom/"),
    @("Bitcoin is flirting with `$20,000 again. How high will it go this time?", "https://mashable.com/article/bitcoin-20000/"),
    @("You can now spend Chinas digital currency at an online store", "https://www.engadget.com/jd-com-supports-china-digital-currency-223718682.html"),
    @("Bitcoin Could Hit `$20K Again But This Time No One Cares", "https://gizmodo.com/bitcoin-could-hit-20k-again-but-this-time-no-one-cares-1845757159"),
    @("Bitcoin hits all-time high rising above `$20,000", "https://www.bbc.co.uk/news/business-55343657"),
    @("With investors expecting a Latin American cryptocurrency boom, Mexico’s Bitso raises `$62 million", "http://techcrunch.com/2020/12/09/with-investors-expecting-a-latin-american-cryptocurrency-boom-mexicos-bitso-raises-62-million/"),
    @("As Crypto comes back, Binance-backed Injective Protocol launches Testnet for its DeFi trading platform", "http://techcrunch.com/2020/12/03/as-crypto-comes-back-binance-backed-injective-protocol-launches-testnet-for-its-defi-trading-platform/"),
    @("Bitcoin soars after BlackRock says it could replace gold", "https://www.cnn.com/2020/11/20/investing/bitcoin-prices-gold-blackrock/index.html"),
    @("Online Shopping, Virus in Winter, Stock Market: Your Monday Evening Briefing", "https://www.nytimes.com/2020/11/30/briefing/online-shopping-virus-in-winter-stock-market.html"),
    @("Months later, the great Twitter hack still boggles my mind", "https://www.theverge.com/22163643/twitter-hack-bitcoin-scam-july-2020-elon-musk"),
    @("Bitcoin price soars above `$18,000, the highest level since 2017", "https://arstechnica.com/tech-policy/2020/11/bitcoin-price-soars-above-18000-the-highest-level-since-2017/"),
    @("Bitcoin price passes `$20,000, setting new record", "https://arstechnica.com/tech-policy/2020/12/bitcoin-price-passes-20000-setting-new-record/"),
    @("Bitcoin breaks above `$20,000 for first time - Reuters", "https://www.reuters.com/article/crypto-currency-idUSKBN28Q1UI"),
    @("Bitcoin breaks above `$20,000 for first time - Reuters India", "https://in.reuters.com/article/crypto-currency-idINKBN28Q1W8"),
    @("Bitcoin hits nearly three-year peak, homes in on record - Reuters.com", "https://uk.reuters.com/article/us-crypto-currencies-idUSKBN27X2PQ"),
    @("Cryptocurrency: Bitcoin hits three-year high as investors jump in - BBC News", "https://www.bbc.com/news/business-54982604"),
    @("How American investors are buying bitcoin - Reuters.com", "https://www.reuters.com/video/watch/idOVD7GKP6Z")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
